$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: update values in B1:E1
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2: clear B2, D2, E2 entirely; update C2 value
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = 5.7271619107786353
$ws.Range("D2").ClearContents()
$ws.Range("E2").ClearContents()

# Row 3: update values
$ws.Range("B3").Value = 6.057254819772349
$ws.Range("C3").Value = 7.0835664376194787
$ws.Range("D3").Value = 9.1862437576153866
$ws.Range("E3").Value = 3.644634541626997

# Update selection to B1:E3
$ws.Range("B1:E3").Select()
